$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.732.81"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -1.05%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.885.66"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -1.36%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'0.9999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.16%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'0.7935"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -5.51%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'241.17"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.47%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.9997"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.01%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  -2.45%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'25.44"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -5.10%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.06996"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -0.76%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.08036"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.01%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.7601"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.67%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'5.290"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +1.03%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = "'1.840.91"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -3.69%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'92.13"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = "'29.709.83"
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = "'13.82"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -3.00%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'5.926"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.62%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'243.13"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -1.02%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.000007667"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -1.48%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = "'0.9999"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.09%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = "'8.164"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +16.32%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'2.124.37"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -1.50%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.22%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'0.1673"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +2.57%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.06%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'163.40"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -3.91%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'18.60"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -2.00%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'2.047"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -1.83%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'1.387"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +1.06%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'1.530"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +0.71%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'4.366"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.32%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.05674"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.60%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'4.045"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -1.62%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'1.259"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -2.23%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.7323"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.60%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.9969"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.16%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'2.614"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -3.73%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.01900"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -1.17%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.76%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.4396"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -1.31%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'72.20"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.75%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'5.814"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -3.35%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.9999"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.08%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.8353"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.90%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'102.60"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +1.28%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'1.022.38"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +3.39%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'1.862"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.36%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'9.850"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +0.67%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'7.413"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -2.79%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'2.019.12"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -2.14%  "
$ws.Range('E51').Style = 'Normal'
